# Insert a new weekly price record before row 439 ("Fruta / hortaliza, semanal").
# All existing rows from 439 downward shift down by one (439->440, ..., 500->501);
# the sheet's used range grows from A1:R500 to A1:R501.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 439:500 down one row, leaving row 439 empty for the new record.
$ws.Rows.Item(439).Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A439").Value = 3
$ws.Range("B439").Value = "Femacal de La Calera"
$ws.Range("C439").Value = "Coquimbo"
$ws.Range("D439").Value = 44984
$ws.Range("E439").Value = 5
$ws.Range("F439").Value = 100112009
$ws.Range("G439").Value = "Acelga"
$ws.Range("H439").Value = "Sin especificar"
$ws.Range("I439").Value = "Primera"
$ws.Range("J439").Value = 260
$ws.Range("K439").Value = 3300
$ws.Range("L439").Value = 3500
$ws.Range("M439").Value = 3392
$ws.Range("N439").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O439").Value = "Provincia de Quillota"
$ws.Range("P439").Value = 565
$ws.Range("Q439").Value = 6
$ws.Range("R439").Value = "Hortaliza"
